$wb = $excel.ActiveWorkbook

# --- Rename the three "Table N - YYYY" sheets (shift the numbering down by two) ---
$wsSheets = $wb.Worksheets
$ws2020 = $wsSheets.Item("Table 5 - 2020")
$ws2020.Name = "Table 3 - 2020"

$ws2021 = $wsSheets.Item("Table 6 - 2021")
$ws2021.Name = "Table 4 - 2021"

$ws2022 = $wsSheets.Item("Table 7 - 2022")
$ws2022.Name = "Table 5 - 2022"

# --- Table 3 - 2020 (formerly "Table 5 - 2020"): selection moved to row 3 ---
$ws2020.Activate()
$ws2020.Rows.Item(3).Select()

# --- Table 4 - 2021 (formerly "Table 6 - 2021"): selection moved to row 3, ---
# --- and row 3 got a slightly shorter custom height ---
$ws2021.Activate()
$ws2021.Rows.Item(3).Select()
$ws2021.Rows.Item(3).RowHeight = 14.25

# --- Table 5 - 2022 (formerly "Table 7 - 2022"): selection moved to row 4 ---
$ws2022.Activate()
$ws2022.Rows.Item(4).Select()

# --- Boxplot: scrolled further down the sheet (selection itself unchanged) ---
$wsBoxplot = $wsSheets.Item("Boxplot")
$wsBoxplot.Activate()
$wsBoxplot.Range("M261").Select()

# --- Table 4_do not use: no longer the active tab, selection moved ---
$wsT4NotUse = $wsSheets.Item("Table 4_do not use")
$wsT4NotUse.Activate()
$wsT4NotUse.Range("T18").Select()

# --- Table 1 ends up as the active tab/sheet ---
$wsTable1 = $wsSheets.Item("Table 1")
$wsTable1.Activate()

Write-Output "done"
